$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking strings
# (e.g. "76.74", "0.0750", "7.60") are not auto-converted to numbers,
# which would silently drop significant trailing zeros / change type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.955.68"
$ws.Range("E2").Value = "  +1.70%  "

$ws.Range("D3").Value = "2.368.51"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "0.689"
$ws.Range("E5").Value = "  +6.05%  "

$ws.Range("D6").Value = "241.64"
$ws.Range("E6").Value = "  +3.04%  "

$ws.Range("D7").Value = "76.74"
$ws.Range("E7").Value = "  +7.35%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "0.633"
$ws.Range("E9").Value = "  +27.54%  "

$ws.Range("E10").Value = "  +5.50%  "

$ws.Range("D11").Value = "57.41"
$ws.Range("E11").Value = "  +1.03%  "

$ws.Range("D12").Value = "33.04"
$ws.Range("E12").Value = "  +21.12%  "

$ws.Range("D13").Value = "7.60"
$ws.Range("E13").Value = "  +19.67%  "

$ws.Range("E14").Value = "  +1.76%  "

$ws.Range("D15").Value = "2.729.10"
$ws.Range("E15").Value = "  +0.97%  "

$ws.Range("D16").Value = "16.92"
$ws.Range("E16").Value = "  +4.12%  "

$ws.Range("D17").Value = "0.926"
$ws.Range("E17").Value = "  +7.20%  "

$ws.Range("D18").Value = "2.374.18"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").Value = "43.952.71"
$ws.Range("E19").Value = "  +1.55%  "

$ws.Range("E20").Value = "  +2.47%  "

$ws.Range("D21").Value = "6.68"
$ws.Range("E21").Value = "  +5.26%  "

$ws.Range("D22").Value = "77.83"
$ws.Range("E22").Value = "  +4.21%  "

$ws.Range("D23").Value = "257.91"
$ws.Range("E23").Value = "  +3.09%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").Value = "2.54"
$ws.Range("E25").Value = "  +3.21%  "

$ws.Range("D26").Value = "11.23"
$ws.Range("E26").Value = "  +11.70%  "

$ws.Range("D27").Value = "3.68"
$ws.Range("E27").Value = "  -3.46%  "

$ws.Range("D28").Value = "1.76"
$ws.Range("E28").Value = "  +14.35%  "

$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  +1.68%  "

$ws.Range("D30").Value = "23.18"
$ws.Range("E30").Value = "  +3.50%  "

$ws.Range("D31").Value = "175.90"
$ws.Range("E31").Value = "  +2.06%  "

$ws.Range("E32").Value = "  -1.65%  "

$ws.Range("E33").Value = "  +5.45%  "

$ws.Range("D34").Value = "5.33"
$ws.Range("E34").Value = "  +6.60%  "

$ws.Range("D35").Value = "0.0750"
$ws.Range("E35").Value = "  +8.55%  "

$ws.Range("D36").Value = "5.37"
$ws.Range("E36").Value = "  +6.15%  "

$ws.Range("D37").Value = "3.83"
$ws.Range("E37").Value = "  +2.63%  "

$ws.Range("D38").Value = "2.46"
$ws.Range("E38").Value = "  +1.23%  "

$ws.Range("D39").Value = "6.51"
$ws.Range("E39").Value = "  -0.40%  "

$ws.Range("D40").Value = "0.0276"
$ws.Range("E40").Value = "  +8.29%  "

$ws.Range("E41").Value = "  +20.15%  "

$ws.Range("D42").Value = "19.05"
$ws.Range("E42").Value = "  +0.23%  "

$ws.Range("D43").Value = "9.06"
$ws.Range("E43").Value = "  +1.52%  "

$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("E45").Value = "  +4.95%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "1.27"
$ws.Range("E46").Value = "  +4.99%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "1.20"
$ws.Range("E47").Value = "  +4.03%  "

$ws.Range("D48").Value = "2.52"
$ws.Range("E48").Value = "  +13.79%  "

$ws.Range("D49").Value = "102.54"
$ws.Range("E49").Value = "  +3.41%  "

$ws.Range("D50").Value = "4.50"
$ws.Range("E50").Value = "  -0.18%  "

$ws.Range("D51").Value = "54.77"
$ws.Range("E51").Value = "  +8.25%  "
